$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 37.16227521392502
$ws.Range("C2").Value = 38.47854333669813
$ws.Range("D2").Value = 35.83546184660323
$ws.Range("E2").Value = 37.16306127174014
$ws.Range("F2").Value = 37.30754659353547
$ws.Range("G2").Value = 37.29334838229598
$ws.Range("H2").Value = 40.34216762033703
$ws.Range("I2").Value = 32.90863680846976
$ws.Range("J2").Value = 37.2851265293913
$ws.Range("K2").Value = 37.05761088918413
$ws.Range("L2").Value = 37.18212647736564
$ws.Range("M2").Value = 37.28271025744524
$ws.Range("N2").Value = 20.32451852977063
$ws.Range("O2").Value = 31.77467961976458
$ws.Range("P2").Value = 41.34592169702232
$ws.Range("Q2").Value = 34.40499174166456
